$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.05 = 7881.45 pesos`n✅ 7881.45 pesos = 2.05 = 929.81 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the N10/O10/N12/O12 rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 488
$wsTasas.Range("O10").Value = 3846.15
$wsTasas.Range("N12").Value = 3850
$wsTasas.Range("O12").Value = 454.201
